$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Join the explorer ... Play for free.</w:t></w:r>
#      </w:p>
# ---------------------------------------------------------------------------

$metaTail = ": Join the explorer on his journey to find lost treasures. Read our review of Adventure Trail, a 5-reel online slot game with 30 active paylines. Play for free."

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

# newly created (still empty) paragraph
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Replace the whole (empty) paragraph's content - including its paragraph
# mark - with the desired run structure via a raw OOXML fragment. Using the
# full (non-collapsed) range here is important: InsertXML *replaces* the
# contents of the exact range it is called on, so calling it on the whole
# paragraph (mark included) cleanly produces the target runs without
# leaking an extra empty run into a neighbouring paragraph.
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:r/>' + `
              '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
              '<w:r><w:t>' + $metaTail + '</w:t></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$metaFullRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$null = $metaFullRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated "Play Adventure Trail Free - Review of Adventure
#    Trail Online Slot" bold paragraph that used to sit right before the
#    closing italic tagline paragraph.
# ---------------------------------------------------------------------------

$oldTitleText = "Play Adventure Trail Free - Review of Adventure Trail Online Slot"
$count = $d.Paragraphs.Count
$dupPara = $d.Paragraphs.Item($count - 1)
if ($dupPara.Range.Text.TrimEnd() -eq $oldTitleText) {
    $dupPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new image
#    prompt, while preserving its italic run formatting.
# ---------------------------------------------------------------------------

$oldTagline = "Join the explorer on his journey to find lost treasures. Read our review of Adventure Trail, a 5-reel online slot game with 30 active paylines. Play for free."
$newTagline = "Create a cartoon-style image featuring a happy Maya warrior with glasses for Adventure Trail online slot game. The Maya warrior should be in traditional clothing, holding a treasure map and standing in the jungle with a big smile on his face. The background should include exotic plants, palm trees and possibly a river. The image should be colorful, with attention-grabbing features to entice potential players to try out the game. The Maya warrior's glasses should also be made prominent to add a unique touch to the image."

$lastCount = $d.Paragraphs.Count
$taglinePara = $d.Paragraphs.Item($lastCount)
if ($taglinePara.Range.Text.TrimEnd() -eq $oldTagline) {
    $taglineRange = $d.Range($taglinePara.Range.Start, $taglinePara.Range.End - 1)
    $taglineRange.Text = $newTagline
}
